$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New custom fill color (FFF828B3) applied to the "Year" header cell (A18),
# mirroring the style used by the other colored label cells.
$ws.Range("A18").Interior.Color = 11741432

# Clear the stray FALSE boolean values left in column F for rows 22,24,26,28,30
$ws.Range("F22").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("F26").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("F30").ClearContents()

# Update the active selection to match the saved view state
$ws.Range("D15").Select()
